{"js": "/*\n * Applies the \"three-digit \u00d7 one-digit\" answer-key update: each table\n * cell holds a single run of text shaped NNN\u00d7N=NNNN; replace the 25\n * updated equations in place using a search/replace per pair so the\n * surrounding run/paragraph formatting (font, size, alignment) is left\n * untouched.\n */\nconst replacements = [\n  [\"407\u00d74=1628\", \"870\u00d75=4350\"],\n  [\"921\u00d77=6447\", \"563\u00d79=5067\"],\n  [\"997\u00d77=6979\", \"945\u00d72=1890\"],\n  [\"253\u00d77=1771\", \"378\u00d74=1512\"],\n  [\"453\u00d72=906\", \"304\u00d77=2128\"],\n  [\"219\u00d79=1971\", \"654\u00d75=3270\"],\n  [\"723\u00d75=3615\", \"824\u00d79=7416\"],\n  [\"668\u00d72=1336\", \"578\u00d75=2890\"],\n  [\"758\u00d78=6064\", \"196\u00d77=1372\"],\n  [\"205\u00d73=615\", \"438\u00d74=1752\"],\n  [\"652\u00d76=3912\", \"921\u00d76=5526\"],\n  [\"255\u00d79=2295\", \"494\u00d78=3952\"],\n  [\"759\u00d75=3795\", \"855\u00d77=5985\"],\n  [\"782\u00d77=5474\", \"729\u00d72=1458\"],\n  [\"506\u00d72=1012\", \"138\u00d75=690\"],\n  [\"556\u00d74=2224\", \"816\u00d73=2448\"],\n  [\"946\u00d78=7568\", \"948\u00d78=7584\"],\n  [\"665\u00d74=2660\", \"261\u00d73=783\"],\n  [\"773\u00d75=3865\", \"423\u00d72=846\"],\n  [\"823\u00d76=4938\", \"343\u00d76=2058\"],\n  [\"509\u00d76=3054\", \"816\u00d72=1632\"],\n  [\"606\u00d75=3030\", \"750\u00d76=4500\"],\n  [\"893\u00d74=3572\", \"599\u00d77=4193\"],\n  [\"744\u00d74=2976\", \"660\u00d73=1980\"],\n  [\"154\u00d74=616\", \"621\u00d79=5589\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Applies the \"three-digit x one-digit\" answer-key update: each table\n# cell holds a single run of text shaped NNNxN=NNNN; replace the 25\n# updated equations in place via Find/Replace so surrounding run and\n# paragraph formatting (font, size, alignment) is left untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"407\u00d74=1628\"; New = \"870\u00d75=4350\" },\n    @{ Old = \"921\u00d77=6447\"; New = \"563\u00d79=5067\" },\n    @{ Old = \"997\u00d77=6979\"; New = \"945\u00d72=1890\" },\n    @{ Old = \"253\u00d77=1771\"; New = \"378\u00d74=1512\" },\n    @{ Old = \"453\u00d72=906\"; New = \"304\u00d77=2128\" },\n    @{ Old = \"219\u00d79=1971\"; New = \"654\u00d75=3270\" },\n    @{ Old = \"723\u00d75=3615\"; New = \"824\u00d79=7416\" },\n    @{ Old = \"668\u00d72=1336\"; New = \"578\u00d75=2890\" },\n    @{ Old = \"758\u00d78=6064\"; New = \"196\u00d77=1372\" },\n    @{ Old = \"205\u00d73=615\"; New = \"438\u00d74=1752\" },\n    @{ Old = \"652\u00d76=3912\"; New = \"921\u00d76=5526\" },\n    @{ Old = \"255\u00d79=2295\"; New = \"494\u00d78=3952\" },\n    @{ Old = \"759\u00d75=3795\"; New = \"855\u00d77=5985\" },\n    @{ Old = \"782\u00d77=5474\"; New = \"729\u00d72=1458\" },\n    @{ Old = \"506\u00d72=1012\"; New = \"138\u00d75=690\" },\n    @{ Old = \"556\u00d74=2224\"; New = \"816\u00d73=2448\" },\n    @{ Old = \"946\u00d78=7568\"; New = \"948\u00d78=7584\" },\n    @{ Old = \"665\u00d74=2660\"; New = \"261\u00d73=783\" },\n    @{ Old = \"773\u00d75=3865\"; New = \"423\u00d72=846\" },\n    @{ Old = \"823\u00d76=4938\"; New = \"343\u00d76=2058\" },\n    @{ Old = \"509\u00d76=3054\"; New = \"816\u00d72=1632\" },\n    @{ Old = \"606\u00d75=3030\"; New = \"750\u00d76=4500\" },\n    @{ Old = \"893\u00d74=3572\"; New = \"599\u00d77=4193\" },\n    @{ Old = \"744\u00d74=2976\"; New = \"660\u00d73=1980\" },\n    @{ Old = \"154\u00d74=616\"; New = \"621\u00d79=5589\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $pair.New\n    $ok = $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n    if (-not $ok) {\n        throw \"Replacement failed for: $($pair.Old)\"\n    }\n}\n"}
